# Daily attendance processing - 2025-10-13 22:26:44
# Re-order "Recorded By" email lists and refresh attendance counts / averages

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (ANATOMY, C1, session 2) - reorder recorders
$ws.Range("G3").Value = "eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 4 (ANATOMY, C1, session 3) - reorder recorders
$ws.Range("G4").Value = "rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 10 - Average Attendance % updated.
# These "xx.x%" labels are stored as plain text in the workbook (not numeric
# percentages). A direct .Value assignment of a percent-looking string gets
# auto-coerced into a numeric percentage by the calc engine, which would
# change the cell's type/style. Route the literal text through a formula and
# then flatten it back to a static value via copy / paste-values so the cell
# keeps its original text type and style.
$ws.Range("L10").Formula = '="44.9%"'
$ws.Range("L10").Copy()
$ws.Range("L10").PasteSpecial(-4163)  # xlPasteValues

# Row 11 (BIOCHEMISTRY LAB/CBL, C1, session 1) - add System recorder, refresh attendance
$ws.Range("G11").Value = "salma.elgendy.std@med.asu.edu.eg, System"
$ws.Range("H11").Value = "159/221"

# Row 12 (HISTOLOGY, C1, session 1) - reorder recorders
$ws.Range("G12").Value = "Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# Row 15 - class statistics attendance % updated (see note on row 10 above)
$ws.Range("S15").Formula = '="52.2%"'
$ws.Range("S15").Copy()
$ws.Range("S15").PasteSpecial(-4163)  # xlPasteValues

# Row 16 - class statistics attendance % updated (see note on row 10 above)
$ws.Range("S16").Formula = '="38.8%"'
$ws.Range("S16").Copy()
$ws.Range("S16").PasteSpecial(-4163)  # xlPasteValues

# Row 19 (PHYSIOLOGY, C1, session 1) - reorder recorders
$ws.Range("G19").Value = "marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# Row 25 (ANATOMY, C2, session 2) - reorder recorders
$ws.Range("G25").Value = "eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 26 (ANATOMY, C2, session 3) - reorder recorders
$ws.Range("G26").Value = "rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 33 (BIOCHEMISTRY LAB/CBL, C2, session 1) - add System recorder, refresh attendance
$ws.Range("G33").Value = "salma.elgendy.std@med.asu.edu.eg, System"
$ws.Range("H33").Value = "146/246"

# Row 34 (HISTOLOGY, C2, session 1) - reorder recorders
$ws.Range("G34").Value = "Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# Row 41 (PHYSIOLOGY, C2, session 1) - reorder recorders
$ws.Range("G41").Value = "yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# Row 42 (PHYSIOLOGY, C2, session 2) - reorder recorders
$ws.Range("G42").Value = "marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
